$d = $word.ActiveDocument

# Donor range: an existing run with the exact target formatting
# (rFonts asciiTheme/hAnsiTheme="majorHAnsi", cs="Calibri Light", sz/szCs=24, rtl, no hint)
# -- the dash right after "KWS " further down in the document has exactly this formatting.
$donor = $d.Content
$donor.Find.Execute("KWS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$donor.Collapse(0)
$donor.MoveEnd(1, 2)
$donor.MoveStart(1, 1)
$ft = $donor.FormattedText

# Target insertion point: right after "נעמה בורמיל – " (name, dash, trailing space run)
$target = $d.Content
$target.Find.Execute("נעמה בורמיל – ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$s = $target.Start

# Paste the donor's formatted (single character) content at the insertion point, then
# overwrite its text with Naama's ID number, keeping the copied run formatting intact.
$target.FormattedText = $ft
$newRun = $d.Range($s, $s + 1)
$newRun.Text = "318319738"
